# Append round-15 match results (rows 140-149) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = @(
    @(138, "Atalanta",       "Milan",    2, 1, 2.11, 0.72, 2.71, 0.96, 0, 0, 0.6,  0.24, 0.84, 3),
    @(139, "Inter",          "Parma",    3, 1, 1.77, 0.38, 2.13, 0.37, 0, 0, 0.36, 0.01, 0.36, 4),
    @(140, "Genoa",          "Torino",   0, 0, 0.49, 0.59, 0.41, 0.57, 0, 0, 0.08, 0.02, 0.1,  0),
    @(141, "Juventus",       "Bologna",  2, 2, 1.02, 0.93, 0.88, 0.59, 0, 0, 0.14, 0.34, 0.48, 4),
    @(142, "Roma",           "Lecce",    4, 1, 1.9,  0.88, 2.52, 0.93, 0, 1, 0.62, 0.05, 0.67, 4),
    @(143, "Fiorentina",     "Cagliari", 1, 0, 0.32, 0.84, 0.45, 1.13, 0, 1, 0.13, 0.29, 0.42, 2),
    @(144, "Hellas Verona",  "Empoli",   1, 4, 0.9399999999999999, 1.31, 0.82, 1.29, 0, 0, 0.12, 0.02, 0.15, 5),
    @(145, "Napoli",         "Lazio",    0, 1, 0.5600000000000001, 0.29, 0.6,  0.37, 0, 0, 0.04, 0.08, 0.12, 1),
    @(146, "Venezia",        "Como",     2, 2, 0.33, 1.05, 0.42, 0.83, 0, 0, 0.09, 0.22, 0.32, 4),
    @(147, "Monza",          "Udinese",  1, 2, 1.32, 1.22, 1.98, 1.3,  0, 0, 0.66, 0.08, 0.75, 3)
)

$startRow = 140
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rec = $data[$i]

    # Column A carries the bordered/bold/centered format used throughout
    # the table; copy it down from the row above (xlPasteFormats = -4122)
    # and then set the value.
    $ws.Cells.Item($row - 1, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($row, 1).Value = $rec[0]

    $ws.Cells.Item($row, 2).Value  = $rec[1]
    $ws.Cells.Item($row, 3).Value  = $rec[2]
    $ws.Cells.Item($row, 4).Value  = $rec[3]
    $ws.Cells.Item($row, 5).Value  = $rec[4]
    $ws.Cells.Item($row, 6).Value  = $rec[5]
    $ws.Cells.Item($row, 7).Value  = $rec[6]
    $ws.Cells.Item($row, 8).Value  = $rec[7]
    $ws.Cells.Item($row, 9).Value  = $rec[8]
    $ws.Cells.Item($row, 10).Value = $rec[9]
    $ws.Cells.Item($row, 11).Value = $rec[10]
    $ws.Cells.Item($row, 12).Value = $rec[11]
    $ws.Cells.Item($row, 13).Value = $rec[12]
    $ws.Cells.Item($row, 14).Value = $rec[13]
    $ws.Cells.Item($row, 15).Value = $rec[14]
}
